# Update "paises.xlsx" (Pais sheet): countries & provincias Spain data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country ranking reshuffle (three-way rotation: Ucrania moves up) ---
$ws.Range("A36").Value = "Ucrania"
$ws.Range("A37").Value = "Bielorrusia"
$ws.Range("A38").Value = "Belgica"

# --- Country ranking reshuffle (Lituania/Estonia swap) ---
$ws.Range("A126").Value = "Lituania"
$ws.Range("A127").Value = "Estonia"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 10:32"

# --- Update numeric case data (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B7").Value = 828990
$ws.Range("C7").Value = 5475
$ws.Range("D7").Value = 620333
$ws.Range("E7").Value = 194984
$ws.Range("G7").Value = 169
$ws.Range("H7").Value = 13673
$ws.Range("B27").Value = 104432
$ws.Range("C27").Value = 2381
$ws.Range("D27").Value = 62138
$ws.Range("E27").Value = 37319
$ws.Range("G27").Value = 74
$ws.Range("H27").Value = 4975
$ws.Range("B36").Value = 67597
$ws.Range("C36").Value = 1022
$ws.Range("D36").Value = 37394
$ws.Range("E36").Value = 28553
$ws.Range("G36").Value = 21
$ws.Range("H36").Value = 1650
$ws.Range("B37").Value = 67366
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 60669
$ws.Range("E37").Value = 6154
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 543
$ws.Range("B38").Value = 66662
$ws.Range("C38").Value = 234
$ws.Range("D38").Value = 17476
$ws.Range("E38").Value = 39353
$ws.Range("G38").Value = 11
$ws.Range("H38").Value = 9833
$ws.Range("B45").Value = 51531
$ws.Range("C45").Value = 334
$ws.Range("E45").Value = 5611
$ws.Range("E52").Value = 3248
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 142
$ws.Range("D63").Value = 16785
$ws.Range("E63").Value = 5979
$ws.Range("G63").Value = 4
$ws.Range("H63").Value = 757
$ws.Range("E116").Value = 1334
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 24
$ws.Range("B124").Value = 2245
$ws.Range("C124").Value = 41
$ws.Range("D124").Value = 1660
$ws.Range("E124").Value = 557
$ws.Range("B126").Value = 2043
$ws.Range("C126").Value = 16
$ws.Range("D126").Value = 1643
$ws.Range("E126").Value = 320
$ws.Range("H126").Value = 80
$ws.Range("B127").Value = 2042
$ws.Range("C127").Value = 4
$ws.Range("D127").Value = 1926
$ws.Range("E127").Value = 47
$ws.Range("H127").Value = 69
$ws.Range("B139").Value = 1224
$ws.Range("C139").Value = 4
$ws.Range("E139").Value = 141
$ws.Range("D162").Value = 441
$ws.Range("E162").Value = 19
